# Atualização automática de GAURAMA.xlsx
$wb = $excel.ActiveWorkbook

# Delete the "Desarquivamentos Pendentes" sheet entirely.
$excel.DisplayAlerts = $false
$ws9 = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$ws9.Delete()

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$ws1 = $wb.Worksheets.Item("Paineis DARQ")
$ws1.Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$ws7 = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$ws7.Name = "RECOLHIMENTO X ELIMINAÇÃO"
